$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 218 (shifts existing rows 218-267 down to 219-268)
$ws.Rows.Item(218).Insert()

# Populate the newly inserted row 218 with the new data record
$ws.Cells.Item(218, 1).Value = 10
$ws.Cells.Item(218, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(218, 3).Value = 'La Araucanía'
$ws.Cells.Item(218, 4).Value = 44642
$ws.Cells.Item(218, 5).Value = 9
$ws.Cells.Item(218, 6).Value = 100112017
$ws.Cells.Item(218, 7).Value = 'Apio'
$ws.Cells.Item(218, 8).Value = 'Americana (o)'
$ws.Cells.Item(218, 9).Value = 'Primera'
$ws.Cells.Item(218, 10).Value = 205
$ws.Cells.Item(218, 11).Value = 9000
$ws.Cells.Item(218, 12).Value = 10000
$ws.Cells.Item(218, 13).Value = 9537
$ws.Cells.Item(218, 14).Value = '$/docena de matas'
$ws.Cells.Item(218, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(218, 16).Value = 1590
$ws.Cells.Item(218, 17).Value = 6
$ws.Cells.Item(218, 18).Value = 'Hortaliza'
